# Auto-generated script applying scraped market-price updates
# to the Sheets workbook (Spriggan_Profits.xlsx) per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 6829.2666
$ws.Range("I9").Value = 8474.166999999999
$ws.Range("K9").Value = 8474.166999999999
$ws.Range("M9").Value = -8305.166999999999
$ws.Range("H29").Value = 2259.8572
$ws.Range("J29").Value = 10999
$ws.Range("L29").Value = 32997
$ws.Range("N29").Value = -33559
$ws.Range("H118").Value = 1198.8889
$ws.Range("I118").Value = 1223.75
$ws.Range("K118").Value = 3671.25
$ws.Range("M118").Value = -2014.25
$ws.Range("H132").Value = 4379.3
$ws.Range("I132").Value = 4379.3
$ws.Range("K132").Value = 13137.9
$ws.Range("M132").Value = -10607.9
$ws.Range("H135").Value = 37038252
$ws.Range("I135").Value = 1211.0952
$ws.Range("K135").Value = 10899.8568
$ws.Range("M135").Value = -8364.8568
$ws.Range("H137").Value = 6722.6665
$ws.Range("I137").Value = 2887
$ws.Range("J137").Value = 7489.8
$ws.Range("K137").Value = 8661
$ws.Range("L137").Value = 22469.4
$ws.Range("M137").Value = -6111
$ws.Range("N137").Value = -27569.4
$ws.Range("H138").Value = 7707.4443
$ws.Range("I138").Value = 5070.6
$ws.Range("J138").Value = 8132.7417
$ws.Range("K138").Value = 15211.8
$ws.Range("L138").Value = 24398.2251
$ws.Range("M138").Value = -10071.8
$ws.Range("N138").Value = -34678.2251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 979986.4399999999
$ws.Range("I2").Value = 1204768
$ws.Range("K2").Value = 1204768
$ws.Range("M2").Value = -1204655
$ws.Range("H32").Value = 9183.666999999999
$ws.Range("I32").Value = 9183.666999999999
$ws.Range("K32").Value = 9183.666999999999
$ws.Range("M32").Value = -8896.666999999999
$ws.Range("H45").Value = 3104.3333
$ws.Range("I45").Value = 3104.3333
$ws.Range("K45").Value = 3104.3333
$ws.Range("M45").Value = -2727.3333
$ws.Range("H110").Value = 1446.8857
$ws.Range("J110").Value = 1436.2693
$ws.Range("L110").Value = 1436.2693
$ws.Range("N110").Value = -5526.2693
$ws.Range("H116").Value = 979986.4399999999
$ws.Range("I116").Value = 1204768
$ws.Range("K116").Value = 1204768
$ws.Range("M116").Value = -1202474
$ws.Range("H122").Value = 1540.5625
$ws.Range("I122").Value = 1517.7858
$ws.Range("K122").Value = 4553.357400000001
$ws.Range("M122").Value = -2103.357400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 979986.4399999999
$ws.Range("I3").Value = 1204768
$ws.Range("K3").Value = 1204768
$ws.Range("M3").Value = -1204654
$ws.Range("H86").Value = 3967.3333
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3967.3333
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 1091.3939
$ws.Range("I107").Value = 1001.93335
$ws.Range("K107").Value = 1001.93335
$ws.Range("M107").Value = 918.06665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 843502.5600000001
$ws.Range("I16").Value = 1685203
$ws.Range("K16").Value = 1685203
$ws.Range("M16").Value = -1684916
$ws.Range("H31").Value = 11344.036
$ws.Range("I31").Value = 7895.25
$ws.Range("K31").Value = 7895.25
$ws.Range("M31").Value = -7600.25
$ws.Range("H34").Value = 11344.036
$ws.Range("I34").Value = 7895.25
$ws.Range("K34").Value = 7895.25
$ws.Range("M34").Value = -7693.25
$ws.Range("H113").Value = 843502.5600000001
$ws.Range("I113").Value = 1685203
$ws.Range("K113").Value = 1685203
$ws.Range("M113").Value = -1683033
$ws.Range("H132").Value = 2467.861
$ws.Range("I132").Value = 2283.2122
$ws.Range("K132").Value = 6849.6366
$ws.Range("M132").Value = -4319.6366
$ws.Range("H134").Value = 6318.5625
$ws.Range("I134").Value = 7232.4614
$ws.Range("K134").Value = 21697.3842
$ws.Range("M134").Value = -19162.3842
$ws.Range("H141").Value = 604548.1
$ws.Range("J141").Value = 649015.2
$ws.Range("L141").Value = 649015.2
$ws.Range("N141").Value = -659375.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4784.8237
$ws.Range("I64").Value = 5188.7144
$ws.Range("J64").Value = 4502.1
$ws.Range("K64").Value = 15566.1432
$ws.Range("L64").Value = 13506.3
$ws.Range("M64").Value = -15296.1432
$ws.Range("N64").Value = -14046.3
$ws.Range("H67").Value = 4784.8237
$ws.Range("I67").Value = 5188.7144
$ws.Range("J67").Value = 4502.1
$ws.Range("K67").Value = 15566.1432
$ws.Range("L67").Value = 13506.3
$ws.Range("M67").Value = -14630.1432
$ws.Range("N67").Value = -15378.3
$ws.Range("H68").Value = 1162.25
$ws.Range("I68").Value = 983
$ws.Range("J68").Value = 1700
$ws.Range("K68").Value = 2949
$ws.Range("L68").Value = 5100
$ws.Range("M68").Value = -2138
$ws.Range("N68").Value = -6722
$ws.Range("H71").Value = 1162.25
$ws.Range("I71").Value = 983
$ws.Range("J71").Value = 1700
$ws.Range("K71").Value = 8847
$ws.Range("L71").Value = 15300
$ws.Range("M71").Value = -4791
$ws.Range("N71").Value = -23412
$ws.Range("H107").Value = 986.56757
$ws.Range("J107").Value = 1181.5333
$ws.Range("L107").Value = 3544.5999
$ws.Range("N107").Value = -7384.5999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 45000
$ws.Range("I26").Value = 45000
$ws.Range("K26").Value = 45000
$ws.Range("M26").Value = -44720
$ws.Range("H50").Value = 45000
$ws.Range("I50").Value = 45000
$ws.Range("K50").Value = 45000
$ws.Range("M50").Value = -44502
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H107").Value = 1127.5
$ws.Range("I107").Value = 353
$ws.Range("K107").Value = 353
$ws.Range("M107").Value = 1567
$ws.Range("H113").Value = 8637.25
$ws.Range("I113").Value = 7160
$ws.Range("K113").Value = 7160
$ws.Range("M113").Value = -4990
$ws.Range("H122").Value = 82176.734
$ws.Range("I122").Value = 94650.08
$ws.Range("K122").Value = 283950.24
$ws.Range("M122").Value = -281500.24
$ws.Range("H126").Value = 10940.417
$ws.Range("I126").Value = 9327.777
$ws.Range("K126").Value = 27983.331
$ws.Range("M126").Value = -25513.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8785
$ws.Range("I7").Value = 7444.25
$ws.Range("K7").Value = 7444.25
$ws.Range("M7").Value = -7332.25
$ws.Range("H16").Value = 2885.8235
$ws.Range("I16").Value = 673.3333
$ws.Range("J16").Value = 5374.875
$ws.Range("K16").Value = 673.3333
$ws.Range("L16").Value = 5374.875
$ws.Range("M16").Value = -503.3333
$ws.Range("N16").Value = -5714.875
$ws.Range("H22").Value = 1661.3871
$ws.Range("I22").Value = 1050.2916
$ws.Range("K22").Value = 1050.2916
$ws.Range("M22").Value = -755.2916
$ws.Range("H27").Value = 1661.3871
$ws.Range("I27").Value = 1050.2916
$ws.Range("K27").Value = 1050.2916
$ws.Range("M27").Value = -943.2916
$ws.Range("H40").Value = 14098.5
$ws.Range("I40").Value = 8969.714
$ws.Range("K40").Value = 8969.714
$ws.Range("M40").Value = -8833.714
$ws.Range("H46").Value = 1698.0834
$ws.Range("I46").Value = 1538.9
$ws.Range("K46").Value = 1538.9
$ws.Range("M46").Value = -1350.9
$ws.Range("H68").Value = 1790050.1
$ws.Range("I68").Value = 2183147.8
$ws.Range("K68").Value = 2183147.8
$ws.Range("M68").Value = -2182398.8
$ws.Range("H71").Value = 1790050.1
$ws.Range("I71").Value = 2183147.8
$ws.Range("K71").Value = 10915739
$ws.Range("M71").Value = -10911995
$ws.Range("H126").Value = 8785
$ws.Range("I126").Value = 7444.25
$ws.Range("K126").Value = 22332.75
$ws.Range("M126").Value = -19862.75
$ws.Range("H132").Value = 51126744
$ws.Range("I132").Value = 57516840
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 172550520
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -172547990
$ws.Range("N132").Value = -23060
$ws.Range("H136").Value = 3201.3635
$ws.Range("I136").Value = 3221.6
$ws.Range("K136").Value = 9664.799999999999
$ws.Range("M136").Value = -7114.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30070
$ws.Range("I54").Value = 30070
$ws.Range("K54").Value = 30070
$ws.Range("M54").Value = -29550
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 1584.8572
$ws.Range("I100").Value = 1481.0385
$ws.Range("K100").Value = 2962.077
$ws.Range("M100").Value = -2421.077
$ws.Range("H136").Value = 2485.9
$ws.Range("J136").Value = 2301.3333
$ws.Range("L136").Value = 6903.999899999999
$ws.Range("N136").Value = -12003.9999
